$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on D-column cells whose new values look numeric,
# so Excel keeps them as text strings (matching the source data) instead of
# silently converting them to numbers.
$dfmt0 = $ws.Range("D4:D6")
$dfmt1 = $ws.Range("D8:D12")
$dfmt2 = $ws.Range("D14:D16")
$dfmt3 = $ws.Range("D18:D20")
$dfmt4 = $ws.Range("D24:D31")
$dfmt5 = $ws.Range("D33:D44")
$dfmt6 = $ws.Range("D46:D48")
$dfmt7 = $ws.Range("D51")
$dfmtUnion = $excel.Union($dfmt0, $dfmt1, $dfmt2, $dfmt3, $dfmt4, $dfmt5, $dfmt6, $dfmt7)
foreach ($area in $dfmtUnion.Areas) {
    $area.NumberFormat = "@"
}

# Apply the updated cell values (price, volume, and the two coin-name/link swaps).
$ws.Range("D2").Value = "30.145.74"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.913.24"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "0.7389"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "244.08"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "0.3128"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "26.90"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Value = "0.06991"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "0.7795"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "0.08000"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.923.49"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "5.296"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "92.36"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "14.43"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "30.154.47"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "5.923"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "242.42"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").Value = "0.000007856"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "2.159.39"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "7.232"
$ws.Range("E24").Value = "  +9.04%  "
$ws.Range("D25").Value = "9.449"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "168.22"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "19.10"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "0.1289"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "2.068"
$ws.Range("E29").Value = "  -5.20%  "
$ws.Range("D30").Value = "1.355"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").Value = "1.548"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "4.110"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "0.05190"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "1.300"
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").Value = "0.7525"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").Value = "2.724"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "0.01946"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "2.802"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "6.384"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "75.26"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4520"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "1.967"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "7.870"
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "0.8396"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.971"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "101.91"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "2.059.15"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "0.1199"
$ws.Range("E51").Value = "  +2.46%  "
